$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("OR-1", 45658, "TEST-ID-15", "SI-1", "PAID"),
    @("OR-2", 45659, "TEST-ID-15", "SI-2", "PAID"),
    @("OR-3", 45660, "TEST-ID-14", "SI-3", "PAID"),
    @("OR-4", 45661, "TEST-ID-15", "SI-4", "PAID"),
    @("OR-5", 45662, "TEST-ID-14", "SI-5", "PAID"),
    @("OR-6", 45663, "TEST-ID-15", "SI-6", "PAID"),
    @("OR-7", 45664, "TEST-ID-15", "SI-7", "PAID"),
    @("OR-8", 45665, "TEST-ID-14", "SI-8", "PAID"),
    @("OR-9", 45666, "TEST-ID-15", "SI-9", "PAID"),
    @("OR-10", 45667, "TEST-ID-15", "SI-10", "PAID"),
    @("OR-11", 45668, "TEST-ID-15", "SI-11", "PAID"),
    @("OR-12", 45669, "TEST-ID-14", "SI-12", "PAID"),
    @("OR-13", 45670, "TEST-ID-15", "SI-13", "PAID"),
    @("OR-14", 45671, "TEST-ID-14", "SI-14", "PAID"),
    @("OR-15", 45672, "TEST-ID-15", "SI-15", "PAID"),
    @("OR-16", 45673, "TEST-ID-14", "SI-16", "PAID"),
    @("OR-17", 45674, "TEST-ID-14", "SI-17", "PAID"),
    @("OR-18", 45675, "TEST-ID-14", "SI-18", "PAID"),
    @("OR-19", 45676, "TEST-ID-14", "SI-19", "PAID"),
    @("OR-20", 45677, "TEST-ID-15", "SI-20", "PAID"),
    @("OR-21", 45678, "TEST-ID-15", "SI-21", "PAID"),
    @("OR-22", 45679, "TEST-ID-15", "SI-22", "PAID"),
    @("OR-23", 45680, "TEST-ID-15", "SI-23", "PAID"),
    @("OR-24", 45681, "TEST-ID-14", "SI-24", "PAID"),
    @("OR-25", 45682, "TEST-ID-15", "SI-25", "PAID"),
    @("OR-26", 45683, "TEST-ID-14", "SI-26", "PAID"),
    @("OR-27", 45684, "TEST-ID-14", "SI-27", "PAID"),
    @("OR-28", 45685, "TEST-ID-15", "SI-28", "PAID"),
    @("OR-29", 45686, "TEST-ID-14", "SI-29", "PAID"),
    @("OR-30", 45687, "TEST-ID-14", "SI-30", "PAID"),
    @("OR-31", 45688, "TEST-ID-14", "SI-31", "PAID"),
    @("OR-32", 45689, "TEST-ID-15", "SI-32", "PAID"),
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

$ws.Range("B2:B33").NumberFormat = "yyyy\-mm\-dd"

$ws.Range("A2").Select() | Out-Null
